$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing formatting / becoming numeric).
$textCells = @("D5", "D6", "D10", "D11", "D12", "D14", "D19", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D31", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Row 2
$ws.Range("D2").Value = "43.188.38"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3
$ws.Range("D3").Value = "2.392.24"
$ws.Range("E3").Value = "  +6.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.43%  "

# Row 5
$ws.Range("D5").Value = "326.93"
$ws.Range("E5").Value = "  +10.05%  "

# Row 6
$ws.Range("D6").Value = "105.41"
$ws.Range("E6").Value = "  -5.87%  "

# Row 7
$ws.Range("E7").Value = "  +2.34%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +6.66%  "

# Row 10
$ws.Range("D10").Value = "41.94"
$ws.Range("E10").Value = "  -4.64%  "

# Row 11
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("D12").Value = "8.70"
$ws.Range("E12").Value = "  -2.23%  "

# Row 13
$ws.Range("E13").Value = "  -0.41%  "

# Row 14
$ws.Range("D14").Value = "17.00"
$ws.Range("E14").Value = "  +11.06%  "

# Row 15
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("D16").Value = "2.751.25"
$ws.Range("E16").Value = "  +6.17%  "

# Row 17
$ws.Range("D17").Value = "2.384.99"
$ws.Range("E17").Value = "  +5.87%  "

# Row 18
$ws.Range("D18").Value = "43.157.45"
$ws.Range("E18").Value = "  +1.00%  "

# Row 19
$ws.Range("D19").Value = "7.71"
$ws.Range("E19").Value = "  +7.52%  "

# Row 20
$ws.Range("E20").Value = "  +2.44%  "

# Row 21
$ws.Range("B21").Value = "PancakeSwap"
$ws.Range("C21").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D21").Value = "3.80"
$ws.Range("E21").Value = "  +8.37%  "

# Row 22
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "76.74"
$ws.Range("E22").Value = "  +1.77%  "

# Row 23
$ws.Range("D23").Value = "274.12"
$ws.Range("E23").Value = "  +6.37%  "

# Row 24
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -1.08%  "

# Row 25
$ws.Range("D25").Value = "9.57"
$ws.Range("E25").Value = "  +7.31%  "

# Row 27
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("D28").Value = "23.13"
$ws.Range("E28").Value = "  +4.04%  "

# Row 29
$ws.Range("D29").Value = "176.20"
$ws.Range("E29").Value = "  +0.54%  "

# Row 30
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("D31").Value = "37.30"
$ws.Range("E31").Value = "  -2.13%  "

# Row 32
$ws.Range("E32").Value = "  +0.60%  "

# Row 33
$ws.Range("D33").Value = "0.0933"
$ws.Range("E33").Value = "  +5.19%  "

# Row 34
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").Value = "  +4.76%  "

# Row 35
$ws.Range("E35").Value = "  +4.55%  "

# Row 36
$ws.Range("D36").Value = "4.89"
$ws.Range("E36").Value = "  -3.49%  "

# Row 37
$ws.Range("D37").Value = "4.16"
$ws.Range("E37").Value = "  -1.39%  "

# Row 38
$ws.Range("D38").Value = "0.0365"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("E39").Value = "  +5.04%  "

# Row 40
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +17.27%  "

# Row 41
$ws.Range("D41").Value = "1.58"
$ws.Range("E41").Value = "  +18.87%  "

# Row 42
$ws.Range("E42").Value = "  +1.16%  "

# Row 43
$ws.Range("D43").Value = "70.07"
$ws.Range("E43").Value = "  -2.78%  "

# Row 44
$ws.Range("D44").Value = "122.02"
$ws.Range("E44").Value = "  +13.20%  "

# Row 45
$ws.Range("E45").Value = "  +0.13%  "

# Row 46
$ws.Range("D46").Value = "12.31"
$ws.Range("E46").Value = "  -1.25%  "

# Row 47
$ws.Range("D47").Value = "89.53"
$ws.Range("E47").Value = "  +41.72%  "

# Row 48
$ws.Range("D48").Value = "9.35"
$ws.Range("E48").Value = "  +7.58%  "

# Row 49
$ws.Range("E49").Value = "  +0.36%  "

# Row 50
$ws.Range("D50").Value = "1.31"
$ws.Range("E50").Value = "  +1.60%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  +3.35%  "
